# COVID-19 Valais figures update (row 374 correction + rows 382-384 newly
# reported daily figures).
#
# Column map: A=Date, B=Cumul cas positifs (formula), C=Nb nouveaux cas
# positifs, D=Nb nouvelles admissions, E=Patients SI, F=Patients intubes,
# G=Patients hosp. hors SI, H=Total hospitalisations (formula),
# I=Nb nouvelles sorties, J=Cumul deces (formula), K=Nb nouveaux deces
# (formula), L=Nb nouveaux deces hopital, M=Nb nouveaux deces extra-hosp.
#
# B, H, J and K are volatile "shared" formulas (they reference TODAY()) so
# writing the plain input cells is enough - the runtime recalculates the
# cached <v> for every dependent cell automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to figures already present in the sheet ---
$ws.Range("C374").Value = 91
$ws.Range("C379").Value = 57
$ws.Range("C380").Value = 78
$ws.Range("C381").Value = 71

# --- Helper: columns L/M carry a Text ("@") number format on these rows,
# even though their stored values must stay numeric (as in the rest of the
# column). Writing .Value straight into such a cell makes the engine store
# a text value instead, so flip the cell to the workbook's plain "Normal"
# style long enough to deposit a real number, then copy the original
# column's number format back from a neighbouring row that already has it
# (this reuses the existing style index instead of minting a new one). ---
function Set-NumericTextCell($ws, [string]$destAddr, [string]$donorAddr, $value) {
    $dest = $ws.Range($destAddr)
    $dest.Style = "Normal"
    $dest.Value = $value
    $ws.Range($donorAddr).Copy()
    $dest.PasteSpecial(-4122)  # xlPasteFormats
}

# --- Row 382 (previously blank placeholder row) ---
$ws.Range("C382").Value = 44
$ws.Range("E382").Value = 5
$ws.Range("F382").Value = 2
$ws.Range("G382").Value = 25
Set-NumericTextCell $ws "L382" "L374" 0
Set-NumericTextCell $ws "M382" "M374" 0

# --- Row 383 ---
$ws.Range("C383").Value = 15
$ws.Range("E383").Value = 4
$ws.Range("F383").Value = 2
$ws.Range("G383").Value = 27
Set-NumericTextCell $ws "L383" "L374" 0
Set-NumericTextCell $ws "M383" "M374" 0

# --- Row 384 ---
$ws.Range("C384").Value = 11
$ws.Range("E384").Value = 6
$ws.Range("F384").Value = 5
$ws.Range("G384").Value = 25
Set-NumericTextCell $ws "L384" "L374" 0
Set-NumericTextCell $ws "M384" "M374" 0

$excel.CutCopyMode = $false
